{"js": "// Edit 1: \"Groupe de 4 \" + bookmark \"_GoBack\" + \"\u00e9tudiants \" -> merge into a\n// single run \"Groupe de 4 \u00e9tudiants \" (bookmark removed, bold kept).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet groupePara = null;\nfor (const p of paras.items) {\n  if (p.text.indexOf(\"Groupe de 4\") !== -1) {\n    groupePara = p;\n    break;\n  }\n}\nif (groupePara) {\n  const groupeRange = groupePara.getRange();\n  groupeRange.insertText(\"Groupe de 4 \\u00e9tudiants \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Edit 2: add a new paragraph right after the paragraph ending in\n// \"...sera consid\u00e9r\u00e9e comme invalide. \" explaining that a frog and a toad\n// must be present on every line, otherwise the configuration is invalid.\n// The new paragraph re-introduces the \"_GoBack\" bookmark right after\n// \"forc\u00e9\" (mirroring where Word leaves it after the last edit).\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\n\nlet invalidePara = null;\nfor (const p of paras2.items) {\n  if (p.text.indexOf(\"consid\\u00e9r\\u00e9e comme invalide\") !== -1) {\n    invalidePara = p;\n    break;\n  }\n}\n\nif (invalidePara) {\n  const newPara = invalidePara.insertParagraph(\n    \"Il faut qu\\u2019il y a forc\\u00e9ment une grenouille et un crapaud sur une ligne sinon la configuration est invalide\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  // Locate \"forc\u00e9\" inside the freshly inserted paragraph and drop a\n  // collapsed \"_GoBack\" bookmark right after it (before \"ment\").\n  const hits = newPara.search(\"forc\\u00e9\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    const afterHit = hits.items[0].insertText(\"\", Word.InsertLocation.after);\n    await context.sync();\n    afterHit.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1 -----------------------------------------------------------\n# \"Groupe de 4 \" + bookmark \"_GoBack\" + \"\u00e9tudiants \" -> merge the two runs\n# into a single run \"Groupe de 4 \u00e9tudiants \" (bold kept, bookmark removed).\n$groupePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Groupe de 4*\") {\n        $groupePara = $p\n        break\n    }\n}\n\nif ($groupePara -ne $null) {\n    $r = $groupePara.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $wasBold = $r.Font.Bold\n    $r.Delete()\n    $r2 = $groupePara.Range\n    $r2.MoveEnd(1, -1) | Out-Null\n    $r2.InsertAfter(\"Groupe de 4 \u00e9tudiants \")\n    $r3 = $groupePara.Range\n    $r3.MoveEnd(1, -1) | Out-Null\n    $r3.Font.Bold = $wasBold\n}\n\n# --- Edit 2 -------------------------------------------------------------\n# Add a new paragraph right after the paragraph ending in\n# \"...sera consid\u00e9r\u00e9e comme invalide. \" and re-create the \"_GoBack\"\n# bookmark right after \"forc\u00e9\" (mirroring where Word leaves it after the\n# last edit made by the author).\n$invalidePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*consid\u00e9r\u00e9e comme invalide*\") {\n        $invalidePara = $p\n        break\n    }\n}\n\nif ($invalidePara -ne $null) {\n    $r = $invalidePara.Range\n    $r.InsertParagraphAfter()\n    $newPara = $invalidePara.Next()\n    $newPara.Range.Text = \"Il faut qu\u2019il y a forc\u00e9ment une grenouille et un crapaud sur une ligne sinon la configuration est invalide\"\n\n    $searchRange = $newPara.Range.Duplicate()\n    $searchRange.Find.Execute(\"forc\u00e9\") | Out-Null\n    $searchRange.Collapse(0) | Out-Null\n    $d.Bookmarks.Add(\"_GoBack\", $searchRange) | Out-Null\n}\n"}
